# Trade #63 closed at 2026-02-17 15:44:36 - unknown UNKNOWN +0.000%
#
# Updates the rolling P&L / trade-count summaries after a new trade closed,
# and appends the new trade's row to both the "All Trades" and
# "MarketMaking" detail logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - headline numbers move by the new trade's P&L.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.38   # Current Capital
$summary.Range("B4").Value = 0.38      # Total P&L $
$summary.Range("B6").Value = 63        # Total Trades
$summary.Range("B7").Value = 21        # Winning Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4) picks up the
#    same deltas as the summary sheet.
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.38     # Capital
$status.Range("D4").Value = 63         # Trades
$status.Range("E4").Value = 0.38       # P&L $
$status.Range("F4").Value = 0.38       # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------
# 3) "All Trades" and "MarketMaking" detail logs - append trade #63 as
#    a new row 64 with identical data in both sheets.
# ---------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 64

    $ws.Range("A$row").Value = 63
    # Leading apostrophe forces this to be stored as literal text instead
    # of being auto-parsed into a date serial number.
    $ws.Range("B$row").Value = "'2026-02-17"
    $ws.Range("C$row").Value = "15:44:29"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("F$row").Value = 0.95
    $ws.Range("G$row").Value = 0.97
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = 2.1053
    $ws.Range("J$row").Value = 0.02
    $ws.Range("K$row").Value = 100.38
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.15
}
